$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115 (pushes the existing row 115 down to 116,
# carrying its formatting/style with it).
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the latest weekly price entry.
$ws.Range("A115").Value = 5
$ws.Range("B115").Value = "Macroferia Regional de Talca"
$ws.Range("C115").Value = "Maule"
$ws.Range("D115").Value = 45121
$ws.Range("E115").Value = 7
$ws.Range("F115").Value = 100112013
$ws.Range("G115").Value = "Alcachofa"
$ws.Range("H115").Value = "Madrigal"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 300
$ws.Range("K115").Value = 15000
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = 15000
$ws.Range("N115").Value = "`$/caja 40 unidades"
$ws.Range("O115").Value = "Provincia del Elquí"
$ws.Range("P115").Value = 375
$ws.Range("Q115").Value = 40
$ws.Range("R115").Value = "Hortaliza"
